$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.041925038812121
$ws.Cells.Item(2, 4).Value2 = 1.050579452143432
$ws.Cells.Item(2, 5).Value2 = 1.055529133488696
$ws.Cells.Item(2, 6).Value2 = 1.062973901854342
$ws.Cells.Item(2, 9).Value2 = 1.042337975978828
$ws.Cells.Item(2, 10).Value2 = 1.047003623203102
$ws.Cells.Item(2, 11).Value2 = 1.05333333687903
$ws.Cells.Item(2, 12).Value2 = 1.058269341153489
$ws.Cells.Item(2, 13).Value2 = 1.065693797644145
$ws.Cells.Item(2, 14).Value2 = 1.019555802400911
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.042874809112019
$ws.Cells.Item(3, 4).Value2 = 1.05133041677942
$ws.Cells.Item(3, 5).Value2 = 1.056426325765546
$ws.Cells.Item(3, 6).Value2 = 1.063864524263542
$ws.Cells.Item(3, 9).Value2 = 1.042544297928471
$ws.Cells.Item(3, 10).Value2 = 1.04759962028164
$ws.Cells.Item(3, 11).Value2 = 1.053896934686771
$ws.Cells.Item(3, 12).Value2 = 1.058979779180804
$ws.Cells.Item(3, 13).Value2 = 1.066399150296461
$ws.Cells.Item(3, 14).Value2 = 1.019754837191647
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.043489812084441
$ws.Cells.Item(4, 4).Value2 = 1.051816380963933
$ws.Cells.Item(4, 5).Value2 = 1.057007552862196
$ws.Cells.Item(4, 6).Value2 = 1.064441286907223
$ws.Cells.Item(4, 9).Value2 = 1.042676199954207
$ws.Cells.Item(4, 10).Value2 = 1.04798508137158
$ws.Cells.Item(4, 11).Value2 = 1.054260994856939
$ws.Cells.Item(4, 12).Value2 = 1.059439527246414
$ws.Cells.Item(4, 13).Value2 = 1.066855392174973
$ws.Cells.Item(4, 14).Value2 = 1.019883519279232
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.043748463203166
$ws.Cells.Item(5, 4).Value2 = 1.052020688408996
$ws.Cells.Item(5, 5).Value2 = 1.057252063020434
$ws.Cells.Item(5, 6).Value2 = 1.064683869073666
$ws.Cells.Item(5, 9).Value2 = 1.042731267328746
$ws.Cells.Item(5, 10).Value2 = 1.048147083007652
$ws.Cells.Item(5, 11).Value2 = 1.054413894935348
$ws.Cells.Item(5, 12).Value2 = 1.059632815449192
$ws.Cells.Item(5, 13).Value2 = 1.067047154685413
$ws.Cells.Item(5, 14).Value2 = 1.01993759130012
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.0437918979258
$ws.Cells.Item(6, 4).Value2 = 1.05205499296123
$ws.Cells.Item(6, 5).Value2 = 1.057293126827965
$ws.Cells.Item(6, 6).Value2 = 1.064724606162241
$ws.Cells.Item(6, 9).Value2 = 1.042740490831066
$ws.Cells.Item(6, 10).Value2 = 1.048174281056847
$ws.Cells.Item(6, 11).Value2 = 1.054439558654715
$ws.Cells.Item(6, 12).Value2 = 1.059665269948012
$ws.Cells.Item(6, 13).Value2 = 1.067079349980409
$ws.Cells.Item(6, 14).Value2 = 1.019946668698907
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.043493267786472
$ws.Cells.Item(7, 4).Value2 = 1.051819110899992
$ws.Cells.Item(7, 5).Value2 = 1.057010819382971
$ws.Cells.Item(7, 6).Value2 = 1.064444527867096
$ws.Cells.Item(7, 9).Value2 = 1.042676937277828
$ws.Cells.Item(7, 10).Value2 = 1.047987246228173
$ws.Cells.Item(7, 11).Value2 = 1.054263038509289
$ws.Cells.Item(7, 12).Value2 = 1.059442109934567
$ws.Cells.Item(7, 13).Value2 = 1.066857954679956
$ws.Cells.Item(7, 14).Value2 = 1.019884241894578
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.042245926818295
$ws.Cells.Item(8, 4).Value2 = 1.050833235043715
$ws.Cells.Item(8, 5).Value2 = 1.055832201418232
$ws.Cells.Item(8, 6).Value2 = 1.063274793592141
$ws.Cells.Item(8, 9).Value2 = 1.042408034748485
$ws.Cells.Item(8, 10).Value2 = 1.047205081719256
$ws.Cells.Item(8, 11).Value2 = 1.053523936196943
$ws.Cells.Item(8, 12).Value2 = 1.058509426542279
$ws.Cells.Item(8, 13).Value2 = 1.065932209005756
$ws.Cells.Item(8, 14).Value2 = 1.019623089013219
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.040051352190976
$ws.Cells.Item(9, 4).Value2 = 1.04909636428316
$ws.Cells.Item(9, 5).Value2 = 1.053760623611509
$ws.Cells.Item(9, 6).Value2 = 1.06121723930214
$ws.Cells.Item(9, 9).Value2 = 1.04192195481812
$ws.Cells.Item(9, 10).Value2 = 1.045825404909673
$ws.Cells.Item(9, 11).Value2 = 1.052216803663048
$ws.Cells.Item(9, 12).Value2 = 1.056866332031981
$ws.Cells.Item(9, 13).Value2 = 1.064299691564024
$ws.Cells.Item(9, 14).Value2 = 1.019162103431637
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.038590646407887
$ws.Cells.Item(10, 4).Value2 = 1.047938782498244
$ws.Cells.Item(10, 5).Value2 = 1.052383215179572
$ws.Cells.Item(10, 6).Value2 = 1.059848087523257
$ws.Cells.Item(10, 9).Value2 = 1.041589713170627
$ws.Cells.Item(10, 10).Value2 = 1.044904734850812
$ws.Cells.Item(10, 11).Value2 = 1.051342260586561
$ws.Cells.Item(10, 12).Value2 = 1.055771281860846
$ws.Cells.Item(10, 13).Value2 = 1.063210589018252
$ws.Cells.Item(10, 14).Value2 = 1.018854263363988
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.037958712449536
$ws.Cells.Item(11, 4).Value2 = 1.047437633045267
$ws.Cells.Item(11, 5).Value2 = 1.051787663178751
$ws.Cells.Item(11, 6).Value2 = 1.059255852484255
$ws.Cells.Item(11, 9).Value2 = 1.041443915053396
$ws.Cells.Item(11, 10).Value2 = 1.044505876412539
$ws.Cells.Item(11, 11).Value2 = 1.050962846916887
$ws.Cells.Item(11, 12).Value2 = 1.055297209314768
$ws.Cells.Item(11, 13).Value2 = 1.062738830961938
$ws.Cells.Item(11, 14).Value2 = 1.018720847313604
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.037724069261163
$ws.Cells.Item(12, 4).Value2 = 1.047251498885273
$ws.Cells.Item(12, 5).Value2 = 1.051566581425472
$ws.Cells.Item(12, 6).Value2 = 1.059035964066918
$ws.Cells.Item(12, 9).Value2 = 1.041389469015646
$ws.Cells.Item(12, 10).Value2 = 1.044357693167843
$ws.Cells.Item(12, 11).Value2 = 1.050821807044695
$ws.Cells.Item(12, 12).Value2 = 1.055121132498914
$ws.Cells.Item(12, 13).Value2 = 1.06256357472681
$ws.Cells.Item(12, 14).Value2 = 1.018671273089198
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.037774397150798
$ws.Cells.Item(13, 4).Value2 = 1.0472914245909
$ws.Cells.Item(13, 5).Value2 = 1.051613998182316
$ws.Cells.Item(13, 6).Value2 = 1.05908312660081
$ws.Cells.Item(13, 9).Value2 = 1.041401160994245
$ws.Cells.Item(13, 10).Value2 = 1.044389480295276
$ws.Cells.Item(13, 11).Value2 = 1.050852065480107
$ws.Cells.Item(13, 12).Value2 = 1.055158900884571
$ws.Cells.Item(13, 13).Value2 = 1.062601168864448
$ws.Cells.Item(13, 14).Value2 = 1.018681907717669
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.037939315026733
$ws.Cells.Item(14, 4).Value2 = 1.047422246822374
$ws.Cells.Item(14, 5).Value2 = 1.051769385775784
$ws.Cells.Item(14, 6).Value2 = 1.05923767451255
$ws.Cells.Item(14, 9).Value2 = 1.041439420441759
$ws.Cells.Item(14, 10).Value2 = 1.044493628134991
$ws.Cells.Item(14, 11).Value2 = 1.050951190729323
$ws.Cells.Item(14, 12).Value2 = 1.055282654441219
$ws.Cells.Item(14, 13).Value2 = 1.062724344715565
$ws.Cells.Item(14, 14).Value2 = 1.018716749848971
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.038040937692416
$ws.Cells.Item(15, 4).Value2 = 1.047502852753561
$ws.Cells.Item(15, 5).Value2 = 1.051865142834432
$ws.Cells.Item(15, 6).Value2 = 1.059332909085187
$ws.Cells.Item(15, 9).Value2 = 1.0414629549245
$ws.Cells.Item(15, 10).Value2 = 1.044557793179937
$ws.Cells.Item(15, 11).Value2 = 1.051012250688757
$ws.Cells.Item(15, 12).Value2 = 1.055358905090725
$ws.Cells.Item(15, 13).Value2 = 1.062800234239396
$ws.Cells.Item(15, 14).Value2 = 1.01873821491938
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.038632597685357
$ws.Cells.Item(16, 4).Value2 = 1.047972044173657
$ws.Cells.Item(16, 5).Value2 = 1.052422758551748
$ws.Cells.Item(16, 6).Value2 = 1.059887405324592
$ws.Cells.Item(16, 9).Value2 = 1.041599348590964
$ws.Cells.Item(16, 10).Value2 = 1.044931201589681
$ws.Cells.Item(16, 11).Value2 = 1.051367425715477
$ws.Cells.Item(16, 12).Value2 = 1.055802746514412
$ws.Cells.Item(16, 13).Value2 = 1.063241894576413
$ws.Cells.Item(16, 14).Value2 = 1.018863115265383
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.039003881269821
$ws.Cells.Item(17, 4).Value2 = 1.048266380916558
$ws.Cells.Item(17, 5).Value2 = 1.052772771342813
$ws.Cells.Item(17, 6).Value2 = 1.060235392261092
$ws.Cells.Item(17, 9).Value2 = 1.041684387047485
$ws.Cells.Item(17, 10).Value2 = 1.045165377328276
$ws.Cells.Item(17, 11).Value2 = 1.051590022722796
$ws.Cells.Item(17, 12).Value2 = 1.05608118176436
$ws.Cells.Item(17, 13).Value2 = 1.063518891894683
$ws.Cells.Item(17, 14).Value2 = 1.018941430272252
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.039220498798825
$ws.Cells.Item(18, 4).Value2 = 1.048438071256705
$ws.Cells.Item(18, 5).Value2 = 1.052977012212214
$ws.Cells.Item(18, 6).Value2 = 1.060438426564744
$ws.Cells.Item(18, 9).Value2 = 1.041733801770699
$ws.Cells.Item(18, 10).Value2 = 1.04530194850098
$ws.Cells.Item(18, 11).Value2 = 1.051719789214789
$ws.Cells.Item(18, 12).Value2 = 1.056243597080194
$ws.Cells.Item(18, 13).Value2 = 1.0636804433144
$ws.Cells.Item(18, 14).Value2 = 1.018987098559149
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.039294368926735
$ws.Cells.Item(19, 4).Value2 = 1.048496614659255
$ws.Cells.Item(19, 5).Value2 = 1.053046667334071
$ws.Cells.Item(19, 6).Value2 = 1.060507666050205
$ws.Cells.Item(19, 9).Value2 = 1.041750619215738
$ws.Cells.Item(19, 10).Value2 = 1.045348512393487
$ws.Cells.Item(19, 11).Value2 = 1.051764024232304
$ws.Cells.Item(19, 12).Value2 = 1.056298977955463
$ws.Cells.Item(19, 13).Value2 = 1.063735525351207
$ws.Cells.Item(19, 14).Value2 = 1.019002668298002
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.038964040453111
$ws.Cells.Item(20, 4).Value2 = 1.048234800460006
$ws.Cells.Item(20, 5).Value2 = 1.052735209550539
$ws.Cells.Item(20, 6).Value2 = 1.060198050400871
$ws.Cells.Item(20, 9).Value2 = 1.041675282539789
$ws.Cells.Item(20, 10).Value2 = 1.045140254505655
$ws.Cells.Item(20, 11).Value2 = 1.051566147463008
$ws.Cells.Item(20, 12).Value2 = 1.056051307384554
$ws.Cells.Item(20, 13).Value2 = 1.063489174394793
$ws.Cells.Item(20, 14).Value2 = 1.01893302900337
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.037890748463044
$ws.Cells.Item(21, 4).Value2 = 1.047383722509104
$ws.Cells.Item(21, 5).Value2 = 1.051723624326504
$ws.Cells.Item(21, 6).Value2 = 1.059192161397278
$ws.Cells.Item(21, 9).Value2 = 1.041428161991883
$ws.Cells.Item(21, 10).Value2 = 1.044462959997906
$ws.Cells.Item(21, 11).Value2 = 1.050922003807493
$ws.Cells.Item(21, 12).Value2 = 1.055246211684867
$ws.Cells.Item(21, 13).Value2 = 1.06268807315855
$ws.Cells.Item(21, 14).Value2 = 1.018706490191688
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.037216420694555
$ws.Cells.Item(22, 4).Value2 = 1.04684870383513
$ws.Cells.Item(22, 5).Value2 = 1.051088369407625
$ws.Cells.Item(22, 6).Value2 = 1.058560263861604
$ws.Cells.Item(22, 9).Value2 = 1.041271109230421
$ws.Cells.Item(22, 10).Value2 = 1.044036948058764
$ws.Cells.Item(22, 11).Value2 = 1.050516376493594
$ws.Cells.Item(22, 12).Value2 = 1.054740102445037
$ws.Cells.Item(22, 13).Value2 = 1.062184249071804
$ws.Cells.Item(22, 14).Value2 = 1.018563954811482
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.037573847525718
$ws.Cells.Item(23, 4).Value2 = 1.047132318575668
$ws.Cells.Item(23, 5).Value2 = 1.051425056737755
$ws.Cells.Item(23, 6).Value2 = 1.058895192596244
$ws.Cells.Item(23, 9).Value2 = 1.041354524740848
$ws.Cells.Item(23, 10).Value2 = 1.044262800855016
$ws.Cells.Item(23, 11).Value2 = 1.050731466446323
$ws.Cells.Item(23, 12).Value2 = 1.05500839197306
$ws.Cells.Item(23, 13).Value2 = 1.062451348626984
$ws.Cells.Item(23, 14).Value2 = 1.018639525009843
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.038982042636147
$ws.Cells.Item(24, 4).Value2 = 1.04824907028221
$ws.Cells.Item(24, 5).Value2 = 1.052752181847191
$ws.Cells.Item(24, 6).Value2 = 1.060214923396696
$ws.Cells.Item(24, 9).Value2 = 1.041679397052049
$ws.Cells.Item(24, 10).Value2 = 1.045151606488348
$ws.Cells.Item(24, 11).Value2 = 1.051576935883361
$ws.Cells.Item(24, 12).Value2 = 1.05606480630393
$ws.Cells.Item(24, 13).Value2 = 1.06350260250468
$ws.Cells.Item(24, 14).Value2 = 1.018936825210811
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.040618293470779
$ws.Cells.Item(25, 4).Value2 = 1.049545334638021
$ws.Cells.Item(25, 5).Value2 = 1.054295539915737
$ws.Cells.Item(25, 6).Value2 = 1.061748723018099
$ws.Cells.Item(25, 9).Value2 = 1.042049064355786
$ws.Cells.Item(25, 10).Value2 = 1.04618224473917
$ws.Cells.Item(25, 11).Value2 = 1.052555283697871
$ws.Cells.Item(25, 12).Value2 = 1.057291055350114
$ws.Cells.Item(25, 13).Value2 = 1.064721875314268
$ws.Cells.Item(25, 14).Value2 = 1.019281371878453
